$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.300877666666666
$ws.Range("H2").Value = 12.902633
$ws.Range("I2").Value = 0.04061703229494078
$ws.Range("J2").Value = 0.04061703229494078
$ws.Range("M2").Value = 18.95273633333333
$ws.Range("N2").Value = 56.858209
$ws.Range("O2").Value = 0.08721078561875104
$ws.Range("P2").Value = 0.08721078561875105
$ws.Range("Q2").Value = 81.51340041825519
$ws.Range("R2").Value = 733.6206037642968
$ws.Range("S2").Value = 0.003542243295943968
$ws.Range("T2").Value = 0.003542243295943969

# Row 3
$ws.Range("G3").Value = 4.300877666666666
$ws.Range("H3").Value = 12.902633
$ws.Range("I3").Value = 0.04061703229494078
$ws.Range("J3").Value = 0.04061703229494078
$ws.Range("O3").Value = 0.04852204497892696
$ws.Range("P3").Value = 0.04852204497892696
$ws.Range("Q3").Value = 45.35215287212666
$ws.Range("R3").Value = 408.16937584914
$ws.Range("S3").Value = 0.001970821467925645
$ws.Range("T3").Value = 0.001970821467925645

# Row 4
$ws.Range("G4").Value = 4.300877666666666
$ws.Range("H4").Value = 12.902633
$ws.Range("I4").Value = 0.04061703229494078
$ws.Range("J4").Value = 0.04061703229494078
$ws.Range("M4").Value = 101.4555613333333
$ws.Range("N4").Value = 304.366684
$ws.Range("O4").Value = 0.4668465309523581
$ws.Range("P4").Value = 0.4668465309523581
$ws.Range("Q4").Value = 436.3479578976635
$ws.Range("R4").Value = 3927.131621078972
$ws.Range("S4").Value = 0.018961920624473
$ws.Range("T4").Value = 0.018961920624473

# Row 5
$ws.Range("G5").Value = 4.300877666666666
$ws.Range("H5").Value = 12.902633
$ws.Range("I5").Value = 0.04061703229494078
$ws.Range("J5").Value = 0.04061703229494078
$ws.Range("M5").Value = 2.410466333333333
$ws.Range("N5").Value = 7.231399000000001
$ws.Range("O5").Value = 0.01109173150200089
$ws.Range("P5").Value = 0.01109173150200089
$ws.Range("Q5").Value = 10.36712081928522
$ws.Range("R5").Value = 93.30408737356699
$ws.Range("S5").Value = 0.0004505132166235823
$ws.Range("T5").Value = 0.0004505132166235824

# Row 6
$ws.Range("G6").Value = 4.300877666666666
$ws.Range("H6").Value = 12.902633
$ws.Range("I6").Value = 0.04061703229494078
$ws.Range("J6").Value = 0.04061703229494078
$ws.Range("M6").Value = 83.95738966666666
$ws.Range("N6").Value = 251.872169
$ws.Range("O6").Value = 0.386328906947963
$ws.Range("P6").Value = 0.386328906947963
$ws.Range("Q6").Value = 361.0904621689973
$ws.Range("R6").Value = 3249.814159520976
$ws.Range("S6").Value = 0.01569153368997458
$ws.Range("T6").Value = 0.01569153368997459

# Row 7
$ws.Range("I7").Value = 0.003347193013613811
$ws.Range("J7").Value = 0.003347193013613811
$ws.Range("M7").Value = 18.95273633333333
$ws.Range("N7").Value = 56.858209
$ws.Range("O7").Value = 0.08721078561875104
$ws.Range("P7").Value = 0.08721078561875105
$ws.Range("Q7").Value = 6.717405703465777
$ws.Range("R7").Value = 60.45665133119199
$ws.Range("S7").Value = 0.0002919113323348553
$ws.Range("T7").Value = 0.0002919113323348553

# Row 8
$ws.Range("I8").Value = 0.003347193013613811
$ws.Range("J8").Value = 0.003347193013613811
$ws.Range("O8").Value = 0.04852204497892696
$ws.Range("P8").Value = 0.04852204497892696
$ws.Range("S8").Value = 0.0001624126499597194
$ws.Range("T8").Value = 0.0001624126499597194

# Row 9
$ws.Range("I9").Value = 0.003347193013613811
$ws.Range("J9").Value = 0.003347193013613811
$ws.Range("M9").Value = 101.4555613333333
$ws.Range("N9").Value = 304.366684
$ws.Range("O9").Value = 0.4668465309523581
$ws.Range("P9").Value = 0.4668465309523581
$ws.Range("Q9").Value = 35.95882696633245
$ws.Range("R9").Value = 323.629442696992
$ws.Range("S9").Value = 0.001562625446833577
$ws.Range("T9").Value = 0.001562625446833577

# Row 10
$ws.Range("I10").Value = 0.003347193013613811
$ws.Range("J10").Value = 0.003347193013613811
$ws.Range("M10").Value = 2.410466333333333
$ws.Range("N10").Value = 7.231399000000001
$ws.Range("O10").Value = 0.01109173150200089
$ws.Range("P10").Value = 0.01109173150200089
$ws.Range("Q10").Value = 0.8543399755457778
$ws.Range("R10").Value = 7.689059779912001
$ws.Range("S10").Value = 0.00003712616619237761
$ws.Range("T10").Value = 0.00003712616619237761

# Row 11
$ws.Range("I11").Value = 0.003347193013613811
$ws.Range("J11").Value = 0.003347193013613811
$ws.Range("M11").Value = 83.95738966666666
$ws.Range("N11").Value = 251.872169
$ws.Range("O11").Value = 0.386328906947963
$ws.Range("P11").Value = 0.386328906947963
$ws.Range("Q11").Value = 29.75696164796355
$ws.Range("R11").Value = 267.812654831672
$ws.Range("S11").Value = 0.001293117418293282
$ws.Range("T11").Value = 0.001293117418293282

# Row 12
$ws.Range("G12").Value = 61.65203333333333
$ws.Range("H12").Value = 184.9561
$ws.Range("I12").Value = 0.582235260574047
$ws.Range("J12").Value = 0.5822352605740471
$ws.Range("M12").Value = 18.95273633333333
$ws.Range("N12").Value = 56.858209
$ws.Range("O12").Value = 0.08721078561875104
$ws.Range("P12").Value = 0.08721078561875105
$ws.Range("Q12").Value = 1168.474732180544
$ws.Range("R12").Value = 10516.2725896249
$ws.Range("S12").Value = 0.05077719448960086
$ws.Range("T12").Value = 0.05077719448960087

# Row 13
$ws.Range("G13").Value = 61.65203333333333
$ws.Range("H13").Value = 184.9561
$ws.Range("I13").Value = 0.582235260574047
$ws.Range("J13").Value = 0.5822352605740471
$ws.Range("O13").Value = 0.04852204497892696
$ws.Range("P13").Value = 0.04852204497892696
$ws.Range("Q13").Value = 650.1120602153333
$ws.Range("R13").Value = 5851.008541937999
$ws.Range("S13").Value = 0.02825124550189117
$ws.Range("T13").Value = 0.02825124550189117

# Row 14
$ws.Range("G14").Value = 61.65203333333333
$ws.Range("H14").Value = 184.9561
$ws.Range("I14").Value = 0.582235260574047
$ws.Range("J14").Value = 0.5822352605740471
$ws.Range("M14").Value = 101.4555613333333
$ws.Range("N14").Value = 304.366684
$ws.Range("O14").Value = 0.4668465309523581
$ws.Range("P14").Value = 0.4668465309523581
$ws.Range("Q14").Value = 6254.941649174711
$ws.Range("R14").Value = 56294.4748425724
$ws.Range("S14").Value = 0.2718145115971361
$ws.Range("T14").Value = 0.2718145115971362

# Row 15
$ws.Range("G15").Value = 61.65203333333333
$ws.Range("H15").Value = 184.9561
$ws.Range("I15").Value = 0.582235260574047
$ws.Range("J15").Value = 0.5822352605740471
$ws.Range("M15").Value = 2.410466333333333
$ws.Range("N15").Value = 7.231399000000001
$ws.Range("O15").Value = 0.01109173150200089
$ws.Range("P15").Value = 0.01109173150200089
$ws.Range("Q15").Value = 148.6101507315444
$ws.Range("R15").Value = 1337.4913565839
$ws.Range("S15").Value = 0.006457997181284856
$ws.Range("T15").Value = 0.006457997181284857

# Row 16
$ws.Range("G16").Value = 61.65203333333333
$ws.Range("H16").Value = 184.9561
$ws.Range("I16").Value = 0.582235260574047
$ws.Range("J16").Value = 0.5822352605740471
$ws.Range("M16").Value = 83.95738966666666
$ws.Range("N16").Value = 251.872169
$ws.Range("O16").Value = 0.386328906947963
$ws.Range("P16").Value = 0.386328906947963
$ws.Range("Q16").Value = 5176.143786308988
$ws.Range("R16").Value = 46585.2940767809
$ws.Range("S16").Value = 0.224934311804134
$ws.Range("T16").Value = 0.224934311804134

# Row 17
$ws.Range("G17").Value = 0.10468
$ws.Range("H17").Value = 0.31404
$ws.Range("I17").Value = 0.0009885868118471018
$ws.Range("J17").Value = 0.0009885868118471018
$ws.Range("M17").Value = 18.95273633333333
$ws.Range("N17").Value = 56.858209
$ws.Range("O17").Value = 0.08721078561875104
$ws.Range("P17").Value = 0.08721078561875105
$ws.Range("Q17").Value = 1.983972439373333
$ws.Range("R17").Value = 17.85575195436
$ws.Range("S17").Value = 0.00008621543251352216
$ws.Range("T17").Value = 0.00008621543251352218

# Row 18
$ws.Range("G18").Value = 0.10468
$ws.Range("H18").Value = 0.31404
$ws.Range("I18").Value = 0.0009885868118471018
$ws.Range("J18").Value = 0.0009885868118471018
$ws.Range("O18").Value = 0.04852204497892696
$ws.Range("P18").Value = 0.04852204497892696
$ws.Range("Q18").Value = 1.1038359448
$ws.Range("R18").Value = 9.9345235032
$ws.Range("S18").Value = 0.00004796825375001908
$ws.Range("T18").Value = 0.00004796825375001908

# Row 19
$ws.Range("G19").Value = 0.10468
$ws.Range("H19").Value = 0.31404
$ws.Range("I19").Value = 0.0009885868118471018
$ws.Range("J19").Value = 0.0009885868118471018
$ws.Range("M19").Value = 101.4555613333333
$ws.Range("N19").Value = 304.366684
$ws.Range("O19").Value = 0.4668465309523581
$ws.Range("P19").Value = 0.4668465309523581
$ws.Range("Q19").Value = 10.62036816037333
$ws.Range("R19").Value = 95.58331344336
$ws.Range("S19").Value = 0.000461518323656071
$ws.Range("T19").Value = 0.000461518323656071

# Row 20
$ws.Range("G20").Value = 0.10468
$ws.Range("H20").Value = 0.31404
$ws.Range("I20").Value = 0.0009885868118471018
$ws.Range("J20").Value = 0.0009885868118471018
$ws.Range("M20").Value = 2.410466333333333
$ws.Range("N20").Value = 7.231399000000001
$ws.Range("O20").Value = 0.01109173150200089
$ws.Range("P20").Value = 0.01109173150200089
$ws.Range("Q20").Value = 0.2523276157733333
$ws.Range("R20").Value = 2.27094854196
$ws.Range("S20").Value = 0.00001096513948342713
$ws.Range("T20").Value = 0.00001096513948342713

# Row 21
$ws.Range("G21").Value = 0.10468
$ws.Range("H21").Value = 0.31404
$ws.Range("I21").Value = 0.0009885868118471018
$ws.Range("J21").Value = 0.0009885868118471018
$ws.Range("M21").Value = 83.95738966666666
$ws.Range("N21").Value = 251.872169
$ws.Range("O21").Value = 0.386328906947963
$ws.Range("P21").Value = 0.386328906947963
$ws.Range("Q21").Value = 8.788659550306665
$ws.Range("R21").Value = 79.09793595276
$ws.Range("S21").Value = 0.0003819196624440624
$ws.Range("T21").Value = 0.0003819196624440624

# Row 22
$ws.Range("G22").Value = 39.47650533333334
$ws.Range("H22").Value = 118.429516
$ws.Range("I22").Value = 0.3728119273055513
$ws.Range("J22").Value = 0.3728119273055513
$ws.Range("M22").Value = 18.95273633333333
$ws.Range("N22").Value = 56.858209
$ws.Range("O22").Value = 0.08721078561875104
$ws.Range("P22").Value = 0.08721078561875105
$ws.Range("Q22").Value = 748.1877969440937
$ws.Range("R22").Value = 6733.690172496844
$ws.Range("S22").Value = 0.03251322106835783
$ws.Range("T22").Value = 0.03251322106835784

# Row 23
$ws.Range("G23").Value = 39.47650533333334
$ws.Range("H23").Value = 118.429516
$ws.Range("I23").Value = 0.3728119273055513
$ws.Range("J23").Value = 0.3728119273055513
$ws.Range("O23").Value = 0.04852204497892696
$ws.Range("P23").Value = 0.04852204497892696
$ws.Range("Q23").Value = 416.2742220292533
$ws.Range("R23").Value = 3746.46799826328
$ws.Range("S23").Value = 0.01808959710540041
$ws.Range("T23").Value = 0.01808959710540041

# Row 24
$ws.Range("G24").Value = 39.47650533333334
$ws.Range("H24").Value = 118.429516
$ws.Range("I24").Value = 0.3728119273055513
$ws.Range("J24").Value = 0.3728119273055513
$ws.Range("M24").Value = 101.4555613333333
$ws.Range("N24").Value = 304.366684
$ws.Range("O24").Value = 0.4668465309523581
$ws.Range("P24").Value = 0.4668465309523581
$ws.Range("Q24").Value = 4005.111008071661
$ws.Range("R24").Value = 36045.99907264495
$ws.Range("S24").Value = 0.1740459549602593
$ws.Range("T24").Value = 0.1740459549602593

# Row 25
$ws.Range("G25").Value = 39.47650533333334
$ws.Range("H25").Value = 118.429516
$ws.Range("I25").Value = 0.3728119273055513
$ws.Range("J25").Value = 0.3728119273055513
$ws.Range("M25").Value = 2.410466333333333
$ws.Range("N25").Value = 7.231399000000001
$ws.Range("O25").Value = 0.01109173150200089
$ws.Range("P25").Value = 0.01109173150200089
$ws.Range("Q25").Value = 95.15678706365378
$ws.Range("R25").Value = 856.4110835728841
$ws.Range("S25").Value = 0.004135129798416651
$ws.Range("T25").Value = 0.004135129798416651

# Row 26
$ws.Range("G26").Value = 39.47650533333334
$ws.Range("H26").Value = 118.429516
$ws.Range("I26").Value = 0.3728119273055513
$ws.Range("J26").Value = 0.3728119273055513
$ws.Range("M26").Value = 83.95738966666666
$ws.Range("N26").Value = 251.872169
$ws.Range("O26").Value = 0.386328906947963
$ws.Range("P26").Value = 0.386328906947963
$ws.Range("Q26").Value = 3314.344340948911
$ws.Range("R26").Value = 29829.0990685402
$ws.Range("S26").Value = 0.1440280243731171
$ws.Range("T26").Value = 0.1440280243731171

Write-Output "Updated 278 cells"
